$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (copy formatting from H1, then set text)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-30
$data = @(
    @(5, 7),
    @(7, 8),
    @(6, 9),
    @(5, 7),
    @(4, 6),
    @(5, 7),
    @(5, 8),
    @(5, 8),
    @(5, 7),
    @(1, 3),
    @(3, 7),
    @(3, 5),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(3, 8),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 5),
    @(1, 1),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(5, 7),
    @(7, 9),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
